$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# G8: "Mandy" -> "Mandya"
$ws.Range("G8").Value = "Mandya"

# Row 23: drop the empty F23 cell (it had no real content)
$ws.Range("F23").ClearContents()

# Row 26: drop the empty F26 cell and correct G26 to just "Mandya"
$ws.Range("F26").ClearContents()
$ws.Range("G26").Value = "Mandya"

# Row 29: drop the empty F29 cell and correct G29 to just "Mandya"
$ws.Range("F29").ClearContents()
$ws.Range("G29").Value = "Mandya"

# G39: "Mandya North" -> "Mandya"
$ws.Range("G39").Value = "Mandya"
